$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "62.288.55"
    "E2" = "  -1.31%  "
    "D3" = "2.442.24"
    "E3" = "  -1.16%  "
    "E4" = "  +0.01%  "
    "D5" = "568.38"
    "E5" = "  -1.32%  "
    "D6" = "145.41"
    "E6" = "  -0.37%  "
    "E7" = "  +0.04%  "
    "D8" = "0.528"
    "E8" = "  -2.52%  "
    "E9" = "  -1.30%  "
    "E10" = "  +0.20%  "
    "D11" = "5.20"
    "E11" = "  -1.62%  "
    "D12" = "0.346"
    "E12" = "  -2.25%  "
    "D13" = "28.56"
    "E13" = "  -1.85%  "
    "D14" = "0.0000173"
    "E14" = "  -3.36%  "
    "D15" = "2.885.09"
    "E15" = "  -1.18%  "
    "D16" = "62.330.32"
    "E16" = "  -1.03%  "
    "D17" = "2.441.09"
    "E17" = "  -1.35%  "
    "D18" = "7.69"
    "E18" = "  -3.18%  "
    "D19" = "10.70"
    "E19" = "  -3.47%  "
    "B20" = "Polkadot"
    "C20" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D20" = "4.11"
    "E20" = "  -0.54%  "
    "B21" = "BitcoinCash"
    "C21" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "D21" = "319.05"
    "E21" = "  -3.30%  "
    "E22" = "  -2.22%  "
    "E23" = "  +0.08%  "
    "D24" = "9.84"
    "E24" = "  +7.31%  "
    "D25" = "64.70"
    "E25" = "  -2.55%  "
    "D26" = "638.35"
    "E26" = "  -3.47%  "
    "D28" = "0.0₃0947"
    "E28" = "  -5.02%  "
    "D29" = "0.997"
    "E29" = "  -0.42%  "
    "D30" = "1.40"
    "E30" = "  -4.77%  "
    "D31" = "7.82"
    "E31" = "  -3.87%  "
    "E32" = "  -3.79%  "
    "E33" = "  -3.75%  "
    "D34" = "0.998"
    "E34" = "  -0.06%  "
    "D35" = "1.48"
    "E35" = "  -3.91%  "
    "D36" = "4.62"
    "E36" = "  -3.60%  "
    "D37" = "150.50"
    "E37" = "  -1.61%  "
    "D38" = "0.364"
    "E38" = "  -2.63%  "
    "D39" = "18.39"
    "E39" = "  -2.24%  "
    "D40" = "5.23"
    "E40" = "  -4.98%  "
    "D41" = "2.68"
    "E41" = "  -1.45%  "
    "D42" = "1.70"
    "E42" = "  -3.68%  "
    "E43" = "  -0.01%  "
    "D44" = "0.0₆0306"
    "E44" = "  +1.99%  "
    "D45" = "151.52"
    "E45" = "  +2.78%  "
    "D46" = "15.32"
    "E46" = "  +1.22%  "
    "D47" = "3.51"
    "E47" = "  -3.18%  "
    "D48" = "0.600"
    "E48" = "  -1.34%  "
    "D49" = "19.97"
    "E49" = "  -4.15%  "
    "D50" = "0.0500"
    "E50" = "  -3.32%  "
    "D51" = "0.0899"
    "E51" = "  -2.42%  "
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
